# Legs and Home Update Sesi 2!
# Update a handful of calibration measurements on Sheet1. Dependent formulas
# (H, I, K, L columns) recalculate automatically when their precedent cells
# are updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# Row 4 (LEG L2)
$ws.Range("D4").Value = 1580
$ws.Range("F4").Value = 1950
$ws.Range("G4").Value = 2050

# Row 5 (LEG L3)
$ws.Range("C5").Value = 1570
$ws.Range("D5").Value = 1550

# Row 6 (LEG R1)
$ws.Range("C6").Value = 1350
$ws.Range("D6").Value = 1680

# Row 7 (LEG R2)
$ws.Range("C7").Value = 1410
$ws.Range("D7").Value = 1380

# Row 8 (LEG R3)
$ws.Range("C8").Value = 1600
$ws.Range("D8").Value = 1250

# Update the saved view/selection to match the new "Home" cell position.
$ws.Range("H12").Select()
